{"js": "// Replace the legacy MERGEFIELD-style field codes (fldChar begin / instrText /\n// fldChar end) that implement the M2Doc \"user content\" zones in the document's\n// default footer with plain literal text runs containing the equivalent\n// M2Doc token syntax: {m:userdoc 'zone1'} and {m:enduserdoc}.\n\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (let s = 0; s < sections.items.length; s++) {\n  const section = sections.items[s];\n\n  // Look at every footer kind (primary/first/even) in case the zone markers\n  // live in a footer other than the default/primary one.\n  const footerTypes = [\"Primary\", \"FirstPage\", \"EvenPages\"];\n\n  for (const footerType of footerTypes) {\n    const footer = section.getFooter(footerType);\n    const paragraphs = footer.paragraphs;\n    paragraphs.load(\"items\");\n    await context.sync();\n\n    // Collect the fields contained in each paragraph of this footer.\n    const fieldsPerParagraph = [];\n    for (let i = 0; i < paragraphs.items.length; i++) {\n      const flds = paragraphs.items[i].getRange().fields;\n      flds.load(\"items\");\n      fieldsPerParagraph.push(flds);\n    }\n    await context.sync();\n\n    // Load the field codes so we know what literal text to substitute.\n    for (const flds of fieldsPerParagraph) {\n      for (let j = 0; j < flds.items.length; j++) {\n        flds.items[j].load(\"code\");\n      }\n    }\n    await context.sync();\n\n    // For every paragraph that is made up of exactly one field (the\n    // begin/instrText/end triplet), replace its whole content with a single\n    // run of literal text \"{<field code>}\".\n    for (let i = 0; i < paragraphs.items.length; i++) {\n      const flds = fieldsPerParagraph[i];\n      if (flds.items.length === 1) {\n        const code = flds.items[0].code.trim();\n        const paragraph = paragraphs.items[i];\n        paragraph.clear();\n        paragraph.insertText(\"{\" + code + \"}\", \"Start\");\n      }\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Replace the legacy MERGEFIELD-style field codes (fldChar begin / instrText /\n# fldChar end) that implement the M2Doc \"user content\" zones in the document's\n# footer(s) with plain literal text runs containing the equivalent M2Doc\n# token syntax: {m:userdoc 'zone1'} and {m:enduserdoc}.\n\n$d = $word.ActiveDocument\n\n# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2, wdHeaderFooterEvenPages = 3\n$footerIndexes = @(1, 2, 3)\n\nfor ($secIdx = 1; $secIdx -le $d.Sections.Count; $secIdx++) {\n    $section = $d.Sections($secIdx)\n\n    foreach ($fIdx in $footerIndexes) {\n        $footer = $section.Footers($fIdx)\n        if (-not $footer.Exists) { continue }\n\n        $paragraphCount = $footer.Range.Paragraphs.Count\n        for ($p = 1; $p -le $paragraphCount; $p++) {\n            $paragraph = $footer.Range.Paragraphs($p)\n            $paragraphRange = $paragraph.Range\n\n            # Only touch paragraphs made up of exactly one field (the\n            # begin/instrText/end triplet) - leave normal text paragraphs\n            # untouched.\n            if ($paragraphRange.Fields.Count -eq 1) {\n                $field = $paragraphRange.Fields(1)\n                $code = $field.Code.Text.Trim()\n\n                # Deleting the field removes the begin/instrText/end runs\n                # entirely, leaving just the paragraph mark behind.\n                $field.Delete()\n\n                # Insert the equivalent M2Doc literal-text token in its place.\n                $paragraphRange.InsertBefore(\"{\" + $code + \"}\")\n            }\n        }\n    }\n}\n"}
